# Add more historic revenue data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in newly available historic JP / Global revenue numbers ---
$ws.Range("B11").Value = 3000000

$ws.Range("B12").Value = 2000000
$ws.Range("C12").Value = 5000000

$ws.Range("B13").Value = 6000000
$ws.Range("C13").Value = 2800000

$ws.Range("B14").Value = 4000000
$ws.Range("C14").Value = 2800000

$ws.Range("B15").Value = 2000000
$ws.Range("C15").Value = 2700000

$ws.Range("B16").Value = 1900000
$ws.Range("C16").Value = 1400000

$ws.Range("B17").Value = 2000000
$ws.Range("C17").Value = 2900000

$ws.Range("B18").Value = 4000000
$ws.Range("C18").Value = 1500000

$ws.Range("B19").Value = 10000000
$ws.Range("C19").Value = 1200000

$ws.Range("B20").Value = 6989822
$ws.Range("C20").Value = 1058556

$ws.Range("B21").Value = 4000000
$ws.Range("C21").Value = 2700000

# Correct a previously mis-entered data point
$ws.Range("B23").Value = 5000000

# --- Move / refresh footnotes that annotate the data source ---
$ws.Range("E22").ClearContents()

$ws.Range("F24").Value = "2022-12 and prior: from u/mee8Ti6Eit"

$ws.Range("E25").ClearContents()
$ws.Range("F25").Value = "2023-01 to 2024-10: from u/numberlockbs"

$ws.Range("E43").ClearContents()
$ws.Range("D43").Value = "* using 2024-07 data from the 2024-08 post"

# --- Restore view state ---
$ws.Range("E46").Select()
$excel.ActiveWindow.ScrollRow = 37
